$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10, shifting the existing
# rows 10-19 down to 11-20 (row 19's data ends up on row 20).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly observation.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44874
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112006
$ws.Range("G10").Value = "Repollo"
$ws.Range("H10").Value = "Copenhague"
$ws.Range("I10").Value = "Tercera"
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 450
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = 475
$ws.Range("N10").Value = "`$/unidad"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 475
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
